$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new institution row (row 40)
$ws.Range("A40").Value = "University of Wisconsin at Madison"
$ws.Range("B40").Value = 240444

# Match formatting of the row above for column A (reuse existing style via copy)
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)

# Column B: start from the same format as the cell above, then shrink the font
# one notch (11pt) - this yields a single new font/style entry instead of one
# per touched property.
$ws.Range("B39").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B40").Font.Size = 11

# Update the view state: scroll so row 20 is the top row, and select C35
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("C35").Select()
